# --- Promote the Sheet2 draft rows (2-8) into Sheet1 as rows 405-411, ---
# --- then repopulate Sheet2 with a new scratch/draft block.            ---

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$sheet1Rows = @(
  ,@('研','nghiên','けんきゅうしゃ','研究者','nhà nghiên cứu',9)
  ,@($null,$null,'けんきゅうしつ','研究室','phòng nghiên cứu',$null)
  ,@('究','cứu','けんきゅうします','研究します','nghiên cứu',7)
  ,@('働','động','はたらきます','働きます','làm việc',13)
  ,@($null,$null,'はたらきすぎ','働きすぎ','làm việc quá nhiều',$null)
  ,@('泳','vịnh','すいえい','水泳','bơi lội',8)
  ,@($null,$null,'およぎます','泳ぎます','bơi(v)',$null)
)

$sheet2Rows = @(
  ,@('部','bộ','ぶちょう','部長','trưởng phòng',10)
  ,@($null,$null,'へや','部屋','căn phòng',$null)
  ,@($null,$null,'ぜんぶ','全部','toàn bộ',$null)
  ,@($null,$null,'いがくぶ','医学部','khoa Y',$null)
  ,@('屋','ốc','へや','部屋','căn phòng',9)
  ,@($null,$null,'～や','～屋','hiệu ~',$null)
  ,@($null,$null,'おくじょう','屋上','tầng thượng',$null)
  ,@('室','thất','きょうしつ','教室','phòng học',9)
  ,@($null,$null,'かいぎしつ','会議室','phòng họp',$null)
  ,@($null,$null,'わしつ','和室','phòng kiểu nhật',$null)
  ,@($null,$null,'けんきゅうしつ','研究室','phòng nghiên cứu',$null)
  ,@($null,$null,'ごうしつ','号室','số phòng',$null)
)

$cols6 = @("A","B","C","D","E","F")

# --- Step 1: Sheet1 gets 7 new rows (405-411), formatted like row 404 ---
for ($i = 0; $i -lt $sheet1Rows.Count; $i++) {
    $r = 405 + $i
    $ws1.Range("A" + ($r - 1) + ":F" + ($r - 1)).Copy($ws1.Range("A" + $r + ":F" + $r))
    $row = $sheet1Rows[$i]
    for ($c = 0; $c -lt 6; $c++) {
        $addr = $cols6[$c] + $r
        if ($row[$c] -eq $null) {
            $ws1.Range($addr).ClearContents()
        } else {
            $ws1.Range($addr).Value = $row[$c]
        }
    }
}

# --- Step 2: Sheet2 is rebuilt in place. Rows 2-9 keep columns A:F;
#             rows 10-27 get the same style/height widened out to column R. ---
for ($r = 3; $r -le 9; $r++) {
    $ws2.Range("A2:F2").Copy($ws2.Range("A" + $r + ":F" + $r))
}
for ($r = 10; $r -le 27; $r++) {
    $ws2.Range("A2:F2").Copy($ws2.Range("A" + $r + ":F" + $r))
    $ws2.Range("A2:F2").Copy($ws2.Range("G" + $r + ":L" + $r))
    $ws2.Range("A2:F2").Copy($ws2.Range("M" + $r + ":R" + $r))
}
$ws2.Range("A2:R27").ClearContents()
$ws2.Range("A2:R27").RowHeight = 18.75

for ($i = 0; $i -lt $sheet2Rows.Count; $i++) {
    $r = 2 + $i
    $row = $sheet2Rows[$i]
    for ($c = 0; $c -lt 6; $c++) {
        $addr = $cols6[$c] + $r
        if ($row[$c] -ne $null) {
            $ws2.Range($addr).Value = $row[$c]
        }
    }
}

# --- Step 3: view state (dimension auto-follows; fix scroll/selection) ---
$ws1.Application.ActiveWindow.ScrollRow = 398
$ws1.Range("E263").Select()
$ws2.Range("F10").Select()

